# Generate Report for Handback
$wb = $excel.ActiveWorkbook

$overviewSheet = $wb.Worksheets.Item("Overview")
$zhSheet = $wb.Worksheets.Item("zh-cn")
$deSheet = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Update the shared "Status" text for the 0eb411dc row (row 3) on every sheet
# that reports it: Overview (B3 + C3), zh-cn (C3), de-de (C3).
$overviewSheet.Range("B3").Value = $newStatus
$overviewSheet.Range("C3").Value = $newStatus
$zhSheet.Range("C3").Value = $newStatus
$deSheet.Range("C3").Value = $newStatus

# Add Error Detail text for zh-cn (row 3, column K)
$zhSheet.Range("K3").Value = "Handback file name: rxgroavp.d1p is different with handoff file name: 0eb411dc-c6b0-4a6a-be20-03ea4d4cf8c7.17041254e66d25123a230ee32d9c546600e54eba.zh-cn."

# Add Error Detail text for de-de (row 3, column K)
$deSheet.Range("K3").Value = "Handback file name: rxgroavp.d1p is different with handoff file name: 0eb411dc-c6b0-4a6a-be20-03ea4d4cf8c7.17041254e66d25123a230ee32d9c546600e54eba.de-de."
